$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2049.2104
$ws.Range("J62").Value = 2042.125
$ws.Range("L62").Value = 2042.125
$ws.Range("N62").Value = -3290.125
$ws.Range("H65").Value = 2049.2104
$ws.Range("J65").Value = 2042.125
$ws.Range("L65").Value = 10210.625
$ws.Range("N65").Value = -16450.625
$ws.Range("H113").Value = 22731456
$ws.Range("I113").Value = 83335180
$ws.Range("K113").Value = 83335180
$ws.Range("M113").Value = -83331926
$ws.Range("H129").Value = 822.03845
$ws.Range("J129").Value = 898.8570999999999
$ws.Range("L129").Value = 2696.5713
$ws.Range("N129").Value = -12696.5713
$ws.Range("H132").Value = 32039.059
$ws.Range("I132").Value = 34784.773
$ws.Range("K132").Value = 104354.319
$ws.Range("M132").Value = -101824.319
$ws.Range("H137").Value = 33027.727
$ws.Range("I137").Value = 3167.1177
$ws.Range("J137").Value = 64754.625
$ws.Range("K137").Value = 9501.3531
$ws.Range("L137").Value = 194263.875
$ws.Range("M137").Value = -6951.3531
$ws.Range("N137").Value = -199363.875
$ws.Range("H138").Value = 2456.742
$ws.Range("J138").Value = 3427.1052
$ws.Range("L138").Value = 10281.3156
$ws.Range("N138").Value = -20561.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 2399
$ws.Range("I21").Value = 497.5
$ws.Range("J21").Value = 3666.6667
$ws.Range("K21").Value = 497.5
$ws.Range("L21").Value = 3666.6667
$ws.Range("M21").Value = -123.5
$ws.Range("N21").Value = -4414.6667
$ws.Range("H74").Value = 2648.926
$ws.Range("I74").Value = 3144.7646
$ws.Range("J74").Value = 1806
$ws.Range("K74").Value = 3144.7646
$ws.Range("L74").Value = 1806
$ws.Range("M74").Value = -2270.7646
$ws.Range("N74").Value = -3554
$ws.Range("H77").Value = 2648.926
$ws.Range("I77").Value = 3144.7646
$ws.Range("J77").Value = 1806
$ws.Range("K77").Value = 15723.823
$ws.Range("L77").Value = 9030
$ws.Range("M77").Value = -11355.823
$ws.Range("N77").Value = -17766
$ws.Range("H132").Value = 24322.092
$ws.Range("I132").Value = 1670.762
$ws.Range("K132").Value = 5012.286
$ws.Range("M132").Value = -2482.286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -717
$ws.Range("H105").Value = 3848130.8
$ws.Range("I105").Value = 2766.6667
$ws.Range("J105").Value = 5001740
$ws.Range("K105").Value = 2766.6667
$ws.Range("L105").Value = 5001740
$ws.Range("M105").Value = -1019.6667
$ws.Range("N105").Value = -5005234
$ws.Range("H134").Value = 40025.406
$ws.Range("I134").Value = 43135.44
$ws.Range("J134").Value = 1150
$ws.Range("K134").Value = 129406.32
$ws.Range("L134").Value = 3450
$ws.Range("M134").Value = -126871.32
$ws.Range("N134").Value = -8520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H31").Value = 9744.223
$ws.Range("I31").Value = 15022.667
$ws.Range("J31").Value = 3146.1667
$ws.Range("K31").Value = 15022.667
$ws.Range("L31").Value = 3146.1667
$ws.Range("M31").Value = -14727.667
$ws.Range("N31").Value = -3736.1667
$ws.Range("H34").Value = 9744.223
$ws.Range("I34").Value = 15022.667
$ws.Range("J34").Value = 3146.1667
$ws.Range("K34").Value = 15022.667
$ws.Range("L34").Value = 3146.1667
$ws.Range("M34").Value = -14820.667
$ws.Range("N34").Value = -3550.1667
$ws.Range("H105").Value = 20834006
$ws.Range("I105").Value = 31250510
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 31250510
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = -31248763
$ws.Range("N105").Value = -4494
$ws.Range("H132").Value = 21371.438
$ws.Range("I132").Value = 23031.393
$ws.Range("J132").Value = 9751.75
$ws.Range("K132").Value = 69094.179
$ws.Range("L132").Value = 29255.25
$ws.Range("M132").Value = -66564.179
$ws.Range("N132").Value = -34315.25
$ws.Range("H134").Value = 7587
$ws.Range("I134").Value = 831.1539
$ws.Range("J134").Value = 51500
$ws.Range("K134").Value = 2493.4617
$ws.Range("L134").Value = 154500
$ws.Range("M134").Value = 41.53830000000016
$ws.Range("N134").Value = -159570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 20
$ws.Range("I7").Value = 20
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 60
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 52
$ws.Range("N7").ClearContents()
$ws.Range("H60").Value = 199.5
$ws.Range("I60").Value = 199.5
$ws.Range("K60").Value = 598.5
$ws.Range("M60").Value = -347.5
$ws.Range("H68").Value = 3978.4
$ws.Range("I68").Value = 1050
$ws.Range("J68").Value = 4155.879
$ws.Range("K68").Value = 3150
$ws.Range("L68").Value = 12467.637
$ws.Range("M68").Value = -2339
$ws.Range("N68").Value = -14089.637
$ws.Range("H71").Value = 3978.4
$ws.Range("I71").Value = 1050
$ws.Range("J71").Value = 4155.879
$ws.Range("K71").Value = 9450
$ws.Range("L71").Value = 37402.911
$ws.Range("M71").Value = -5394
$ws.Range("N71").Value = -45514.911
$ws.Range("H76").Value = 4924.091
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4924.091
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 14772.273
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -15538.273
$ws.Range("H79").Value = 4924.091
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4924.091
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 14772.273
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -17424.273
$ws.Range("H107").Value = 4986.6
$ws.Range("J107").Value = 1062.2273
$ws.Range("L107").Value = 3186.6819
$ws.Range("N107").Value = -7026.6819
$ws.Range("H131").Value = 127426.32
$ws.Range("I131").Value = 794.2857
$ws.Range("J131").Value = 139737.77
$ws.Range("K131").Value = 2382.8571
$ws.Range("L131").Value = 419213.3099999999
$ws.Range("M131").Value = 2657.1429
$ws.Range("N131").Value = -429293.3099999999
$ws.Range("H140").Value = 1806.2222
$ws.Range("I140").Value = 1438.9333
$ws.Range("J140").Value = 3642.6667
$ws.Range("K140").Value = 4316.7999
$ws.Range("L140").Value = 10928.0001
$ws.Range("M140").Value = 863.2001
$ws.Range("N140").Value = -21288.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4000
$ws.Range("J18").Value = 4000
$ws.Range("L18").Value = 4000
$ws.Range("N18").Value = -4586
$ws.Range("H80").Value = 9151.789000000001
$ws.Range("I80").Value = 19534.166
$ws.Range("K80").Value = 19534.166
$ws.Range("M80").Value = -18536.166
$ws.Range("H83").Value = 9151.789000000001
$ws.Range("I83").Value = 19534.166
$ws.Range("K83").Value = 97670.83
$ws.Range("M83").Value = -92678.83
$ws.Range("H95").Value = 24500
$ws.Range("J95").Value = 24500
$ws.Range("L95").Value = 24500
$ws.Range("N95").Value = -29992
$ws.Range("H97").Value = 1916.1904
$ws.Range("I97").Value = 1040.8823
$ws.Range("J97").Value = 5636.25
$ws.Range("K97").Value = 1040.8823
$ws.Range("L97").Value = 5636.25
$ws.Range("M97").Value = -544.8823
$ws.Range("N97").Value = -6628.25
$ws.Range("H122").Value = 4165.3
$ws.Range("J122").Value = 8500
$ws.Range("L122").Value = 25500
$ws.Range("N122").Value = -30400
$ws.Range("H132").Value = 87432.836
$ws.Range("I132").Value = 105030.1
$ws.Range("J132").Value = 65436.25
$ws.Range("K132").Value = 315090.3
$ws.Range("L132").Value = 196308.75
$ws.Range("M132").Value = -312560.3
$ws.Range("N132").Value = -201368.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5197.7617
$ws.Range("J7").Value = 4549.75
$ws.Range("L7").Value = 4549.75
$ws.Range("N7").Value = -4773.75
$ws.Range("H25").Value = 4000
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H101").Value = 11787.333
$ws.Range("J101").Value = 11787.333
$ws.Range("L101").Value = 11787.333
$ws.Range("N101").Value = -18277.333
$ws.Range("H126").Value = 5197.7617
$ws.Range("J126").Value = 4549.75
$ws.Range("L126").Value = 13649.25
$ws.Range("N126").Value = -18589.25
$ws.Range("H132").Value = 1592.8837
$ws.Range("I132").Value = 1175.6875
$ws.Range("J132").Value = 2806.5454
$ws.Range("K132").Value = 3527.0625
$ws.Range("L132").Value = 8419.636200000001
$ws.Range("M132").Value = -997.0625
$ws.Range("N132").Value = -13479.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 4200
$ws.Range("I34").Value = 3000
$ws.Range("J34").Value = 4800
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 4800
$ws.Range("M34").Value = -2797
$ws.Range("N34").Value = -5206
$ws.Range("H107").Value = 1436.7368
$ws.Range("I107").Value = 415.1
$ws.Range("J107").Value = 2571.889
$ws.Range("K107").Value = 1245.3
$ws.Range("L107").Value = 7715.667
$ws.Range("M107").Value = 674.6999999999998
$ws.Range("N107").Value = -11555.667
$ws.Range("H111").Value = 35644
$ws.Range("J111").Value = 35644
$ws.Range("L111").Value = 35644
$ws.Range("N111").Value = -43824
$ws.Range("H122").Value = 1564.9166
$ws.Range("I122").Value = 1410.1904
$ws.Range("K122").Value = 4230.5712
$ws.Range("M122").Value = -1780.5712
$ws.Range("H126").Value = 2249.3333
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -13190
$ws.Range("H132").Value = 3599.75
$ws.Range("I132").Value = 3133.3333
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9399.999899999999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6869.999899999999
$ws.Range("N132").Value = -20057
